# Monopoly map edit: delete unimplemented event/tile-title strings and
# recategorize every tile's eventId/tileTitle into the new simplified
# bucket scheme (EARN / TAX / BANK / SLOT), per the commit
# "edit monopoly map, delete not unimplemented events".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("MapConfig")
$ws2 = $wb.Worksheets.Item("MapEvent")

# --- MapEvent sheet: event 103's description changes from "jail" to
#     "Slot mini game" (the other event rows 100/101/102/104/105 stay put).
$ws2.Range("C7").Value = "Slot mini game"

# --- MapConfig sheet: for every tile row (4..43) set the eventId (col F)
#     and the tileTitle (col I) to their new values. Rows 13 and 23
#     previously had no eventId at all; they now get one too.
$tileData = @{
    4  = @{ F = 100; I = "EARN" }
    5  = @{ F = 100; I = "EARN" }
    6  = @{ F = 100; I = "EARN" }
    7  = @{ F = 101; I = "TAX"  }
    8  = @{ F = 102; I = "BANK" }
    9  = @{ F = 101; I = "TAX"  }
    10 = @{ F = 100; I = "EARN" }
    11 = @{ F = 103; I = "SLOT" }
    12 = @{ F = 100; I = "EARN" }
    13 = @{ F = 100; I = "EARN" }
    14 = @{ F = 100; I = "EARN" }
    15 = @{ F = 100; I = "EARN" }
    16 = @{ F = 100; I = "EARN" }
    17 = @{ F = 101; I = "TAX"  }
    18 = @{ F = 102; I = "BANK" }
    19 = @{ F = 101; I = "TAX"  }
    20 = @{ F = 100; I = "EARN" }
    21 = @{ F = 103; I = "SLOT" }
    22 = @{ F = 100; I = "EARN" }
    23 = @{ F = 100; I = "EARN" }
    24 = @{ F = 100; I = "EARN" }
    25 = @{ F = 100; I = "EARN" }
    26 = @{ F = 100; I = "EARN" }
    27 = @{ F = 100; I = "EARN" }
    28 = @{ F = 102; I = "BANK" }
    29 = @{ F = 101; I = "TAX"  }
    30 = @{ F = 101; I = "TAX"  }
    31 = @{ F = 100; I = "EARN" }
    32 = @{ F = 103; I = "SLOT" }
    33 = @{ F = 100; I = "EARN" }
    34 = @{ F = 100; I = "EARN" }
    35 = @{ F = 100; I = "EARN" }
    36 = @{ F = 100; I = "EARN" }
    37 = @{ F = 100; I = "EARN" }
    38 = @{ F = 102; I = "BANK" }
    39 = @{ F = 101; I = "TAX"  }
    40 = @{ F = 101; I = "TAX"  }
    41 = @{ F = 103; I = "SLOT" }
    42 = @{ F = 100; I = "EARN" }
    43 = @{ F = 100; I = "EARN" }
}

foreach ($row in $tileData.Keys) {
    $entry = $tileData[$row]
    $ws1.Cells.Item($row, 6).Value = $entry.F   # column F = eventId
    $ws1.Cells.Item($row, 9).Value = $entry.I   # column I = tileTitle
}

# Row 43's event count (col G) drops from 10 to 1.
$ws1.Range("G43").Value = 1

# --- Restore the cell selections shown in each sheet's view, and make
#     sure MapConfig (the tab that was selected before) stays selected.
$ws2.Activate()
$ws2.Range("C7").Select() | Out-Null

$ws1.Activate()
$ws1.Range("J20").Select() | Out-Null
